$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Containers")

# Update container name in B9 from "messagestodeveloper" to "messagesubmissions"
$ws.Range("B9").Value = "messagesubmissions"

# Update the active selection to match the authored state
$ws.Range("B10").Select()
